# Archivo.xlsx update ("subiendo archivos en github")
#
# 1) The shared-string record holding the form data has the email changed
#    from deisy.Saenz3@gmail.com to deisy.Saenz5@gmail.com.
# 2) The RegistroUsuario sheet view picks up an active selection on E19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$current = $ws.Range("A1").Value()
$updated = $current.Replace("deisy.Saenz3@gmail.com", "deisy.Saenz5@gmail.com")
$ws.Range("A1").Value = $updated

$ws.Range("E19").Select()
